# "this added last report 29-03-25"
# Update the requisition quantities on Sheet1 and refresh the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (MMST 470): quantity 25 -> 15
$ws.Range("C9").Value = 15

# Row 31 (Router): quantity 5000 -> 6500
$ws.Range("C31").Value = 6500

# Row 32 (Pocket Router): quantity 500 -> (cleared)
$ws.Range("C32").ClearContents()

# Row 43 (CREDIT Lifting): quantity 315964 -> 266862
$ws.Range("C43").Value = 266862

# Move the active selection to C32, matching the saved view state.
$ws.Activate()
$ws.Range("C32").Select()
